# Applies the stock-report correction described in the commit diff:
# - adjusts Qty (F) / Value (G) on ~60 line items (mostly -1 unit, a few larger deltas)
# - swaps a handful of duplicate-named rows back into their original order
#   (B/C/E/F/G got transposed between adjacent rows in the source data)
# - recalculates every "Sub Total:" / "Grand Total:" row (column B) affected by the above
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("F28").Value = 47
$ws.Range("G28").Value = 1444.78

# Row 30
$ws.Range("F30").Value = 134
$ws.Range("G30").Value = 3776.12

# Row 34
$ws.Range("B34").Value = 58196.81

# Row 48
$ws.Range("F48").Value = 233
$ws.Range("G48").Value = 13071.3

# Row 52
$ws.Range("F52").Value = 47
$ws.Range("G52").Value = 2773

# Row 55
$ws.Range("F55").Value = 121
$ws.Range("G55").Value = 6746.96

# Row 61
$ws.Range("F61").Value = 223
$ws.Range("G61").Value = 58142.79

# Row 66
$ws.Range("B66").Value = 201711.94

# Row 161
$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644

# Row 162
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88

# Row 173
$ws.Range("F173").Value = 46
$ws.Range("G173").Value = 3616.06

# Row 193
$ws.Range("B193").Value = 63869.45

# Row 213
$ws.Range("F213").Value = 208
$ws.Range("G213").Value = 26349.44

# Row 216
$ws.Range("F216").Value = 75
$ws.Range("G216").Value = 5572.5

# Row 217
$ws.Range("F217").Value = 33
$ws.Range("G217").Value = 2451.9

# Row 218
$ws.Range("B218").Value = 75508.12

# Row 222
$ws.Range("F222").Value = 778
$ws.Range("G222").Value = 14393

# Row 227
$ws.Range("F227").Value = 30
$ws.Range("G227").Value = 3438

# Row 229
$ws.Range("B229").Value = 24051.63

# Row 262
$ws.Range("F262").Value = 73
$ws.Range("G262").Value = 5718.09

# Row 277
$ws.Range("F277").Value = 7
$ws.Range("G277").Value = 352.87

# Row 288
$ws.Range("F288").Value = 3
$ws.Range("G288").Value = 1665.09

# Row 292
$ws.Range("B292").Value = 64985
$ws.Range("C292").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F292").Value = 12
$ws.Range("G292").Value = 1052.4

# Row 293
$ws.Range("B293").Value = 66196
$ws.Range("C293").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F293").Value = 6
$ws.Range("G293").Value = 526.2

# Row 295
$ws.Range("B295").Value = 114366.05

# Row 297
$ws.Range("B297").Value = 63565
$ws.Range("E297").Value = 109.19
$ws.Range("F297").Value = 60
$ws.Range("G297").Value = 6162.6

# Row 298
$ws.Range("B298").Value = 61610
$ws.Range("E298").Value = 122.71
$ws.Range("F298").Value = -58
$ws.Range("G298").Value = -5957.18

# Row 306
$ws.Range("B306").Value = 57802
$ws.Range("E306").Value = 162.71
$ws.Range("F306").Value = -79
$ws.Range("G306").Value = -11334.92

# Row 307
$ws.Range("B307").Value = 63531
$ws.Range("E307").Value = 152.53
$ws.Range("F307").Value = 27
$ws.Range("G307").Value = 3873.96

# Row 324
$ws.Range("F324").Value = 25
$ws.Range("G324").Value = 4283.25

# Row 325
$ws.Range("F325").Value = 25
$ws.Range("G325").Value = 3779.25

# Row 328
$ws.Range("B328").Value = -12049.27

# Row 359
$ws.Range("F359").Value = 53
$ws.Range("G359").Value = 12725.3

# Row 363
$ws.Range("B363").Value = 74196.50999999999

# Row 367
$ws.Range("F367").Value = 128
$ws.Range("G367").Value = 7769.6

# Row 370
$ws.Range("F370").Value = 207
$ws.Range("G370").Value = 34359.93

# Row 372
$ws.Range("B372").Value = 58348.59

# Row 380
$ws.Range("F380").Value = 37
$ws.Range("G380").Value = 2032.04

# Row 389
$ws.Range("B389").Value = 57302.46

# Row 398
$ws.Range("F398").Value = 90
$ws.Range("G398").Value = 5051.7

# Row 408
$ws.Range("F408").Value = 204
$ws.Range("G408").Value = 3233.4

# Row 410
$ws.Range("F410").Value = 72
$ws.Range("G410").Value = 17357.76

# Row 416
$ws.Range("F416").Value = 72
$ws.Range("G416").Value = 2113.92

# Row 417
$ws.Range("B417").Value = 170667.34

# Row 423
$ws.Range("F423").Value = 54
$ws.Range("G423").Value = 4513.86

# Row 427
$ws.Range("B427").Value = 23426.12

# Row 432
$ws.Range("F432").Value = 91
$ws.Range("G432").Value = 4405.31

# Row 438
$ws.Range("B438").Value = 24638.33

# Row 479
$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52

# Row 480
$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0

# Row 511
$ws.Range("F511").Value = 226
$ws.Range("G511").Value = 22570.62

# Row 512
$ws.Range("F512").Value = 20
$ws.Range("G512").Value = 2371.6

# Row 525
$ws.Range("B525").Value = 119318.83

# Row 527
$ws.Range("F527").Value = 48
$ws.Range("G527").Value = 1589.28

# Row 529
$ws.Range("F529").Value = 118
$ws.Range("G529").Value = 3906.98

# Row 534
$ws.Range("F534").Value = 123
$ws.Range("G534").Value = 5382.48

# Row 535
$ws.Range("B535").Value = 23372.06

# Row 565
$ws.Range("F565").Value = 15
$ws.Range("G565").Value = 4214.25

# Row 573
$ws.Range("B573").Value = 22956.7

# Row 617
$ws.Range("F617").Value = 21
$ws.Range("G617").Value = 1010.52

# Row 620
$ws.Range("F620").Value = 358
$ws.Range("G620").Value = 28135.22

# Row 625
$ws.Range("F625").Value = 320
$ws.Range("G625").Value = 11785.6

# Row 628
$ws.Range("B628").Value = 208962.88

# Row 659
$ws.Range("F659").Value = 38
$ws.Range("G659").Value = 2034.52

# Row 660
$ws.Range("F660").Value = 50
$ws.Range("G660").Value = 1487

# Row 662
$ws.Range("F662").Value = 38
$ws.Range("G662").Value = 3051.78

# Row 668
$ws.Range("B668").Value = 11828.22

# Row 674
$ws.Range("F674").Value = 740
$ws.Range("G674").Value = 120701.4

# Row 680
$ws.Range("B680").Value = 121713.95

# Row 694
$ws.Range("F694").Value = 4
$ws.Range("G694").Value = 1619.08

# Row 706
$ws.Range("F706").Value = 121
$ws.Range("G706").Value = 4741.99

# Row 711
$ws.Range("F711").Value = 12
$ws.Range("G711").Value = 6406.2

# Row 713
$ws.Range("B713").Value = 64808.36

# Row 718
$ws.Range("B718").Value = 2632998.53

# Row 719
$ws.Range("B719").Value = 2632998.53
